# Insert a new weekly price-report row at row 569 of the (single) sheet,
# pushing the existing rows 569-629 down to 570-630, and populate the new
# row with its own data (dimension grows from A1:R629 to A1:R630).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Inserting a whole row shifts everything below it down by one and copies
# the formatting (incl. the date-format style on column D) from the row
# that used to occupy this position - matching how Excel's own
# Rows.Insert works when a user inserts a row above the current one.
$ws.Rows.Item(569).Insert()

# Fill in the new row with the new observation's data.
$ws.Cells.Item(569, 1).Value = 3
$ws.Cells.Item(569, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(569, 3).Value = "Coquimbo"
$ws.Cells.Item(569, 4).Value = 45194
$ws.Cells.Item(569, 5).Value = 5
$ws.Cells.Item(569, 6).Value = 100114013
$ws.Cells.Item(569, 7).Value = "Zanahoria"
$ws.Cells.Item(569, 8).Value = "Sin especificar"
$ws.Cells.Item(569, 9).Value = "Primera"
$ws.Cells.Item(569, 10).Value = 270
$ws.Cells.Item(569, 11).Value = 9500
$ws.Cells.Item(569, 12).Value = 10000
$ws.Cells.Item(569, 13).Value = 9796
$ws.Cells.Item(569, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(569, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(569, 16).Value = 490
$ws.Cells.Item(569, 17).Value = 20
$ws.Cells.Item(569, 18).Value = "Hortaliza"
